# Apply "data up to 12" update to the community-state survey sheet.
# - Corrects several recomputed estimates in rows 119-122
# - Adds new daily rows 123-128 (full state-level data) and 129-130 (date-only rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrected estimates for existing rows (119-122) ---
$ws.Range("M119").Value = 16.3863843
$ws.Range("M120").Value = 16.5547062
$ws.Range("M121").Value = 16.5490652
$ws.Range("F122").Value = 14.899868
$ws.Range("G122").Value = 11.3613703
$ws.Range("H122").Value = 15.8248666
$ws.Range("L122").Value = 11.7010078
$ws.Range("M122").Value = 16.2632995
$ws.Range("T122").Value = 14.8334356
$ws.Range("AN122").Value = 14.6682243
$ws.Range("AV122").Value = 13.9523877
$ws.Range("AW122").Value = 14.6713623
$ws.Range("BB122").Value = 13.4325572

# --- New date labels (column A) for rows 124-130 ---
$ws.Range("A124").Value = "02 06 2020"
$ws.Range("A125").Value = "03 06 2020"
$ws.Range("A126").Value = "04 06 2020"
$ws.Range("A127").Value = "05 06 2020"
$ws.Range("A128").Value = "06 06 2020"
$ws.Range("A129").Value = "07 06 2020"
$ws.Range("A130").Value = "08 06 2020"

# --- New survey data for rows 123-128 ---
# Row 123
$ws.Range("B123").Value = 12.5605164
$ws.Range("C123").Value = 18.7590793
$ws.Range("D123").Value = 16.7497956
$ws.Range("F123").Value = 14.9050471
$ws.Range("G123").Value = 11.3592035
$ws.Range("H123").Value = 15.5485981
$ws.Range("I123").Value = 15.9224559
$ws.Range("J123").Value = 15.546697
$ws.Range("K123").Value = 14.5731478
$ws.Range("L123").Value = 11.8885942
$ws.Range("M123").Value = 16.1794707
$ws.Range("O123").Value = 8.8235294
$ws.Range("P123").Value = 21.6181246
$ws.Range("Q123").Value = 13.4141061
$ws.Range("R123").Value = 16.1735404
$ws.Range("S123").Value = 18.9726432
$ws.Range("T123").Value = 14.6231776
$ws.Range("U123").Value = 15.4351575
$ws.Range("V123").Value = 19.0531787
$ws.Range("W123").Value = 17.3041009
$ws.Range("X123").Value = 15.7569205
$ws.Range("Y123").Value = 10.0152572
$ws.Range("Z123").Value = 14.1996381
$ws.Range("AA123").Value = 17.5769097
$ws.Range("AB123").Value = 14.81863
$ws.Range("AD123").Value = 19.7109308
$ws.Range("AE123").Value = 9.2229922
$ws.Range("AF123").Value = 13.9046305
$ws.Range("AG123").Value = 17.8383689
$ws.Range("AH123").Value = 21.426789
$ws.Range("AI123").Value = 13.4365526
$ws.Range("AJ123").Value = 17.5675617
$ws.Range("AK123").Value = 14.1010263
$ws.Range("AL123").Value = 12.1167964
$ws.Range("AM123").Value = 15.5392094
$ws.Range("AN123").Value = 14.3210989
$ws.Range("AO123").Value = 14.1231491
$ws.Range("AP123").Value = 11.5758286
$ws.Range("AQ123").Value = 12.3378264
$ws.Range("AR123").Value = 7.9882485
$ws.Range("AS123").Value = 15.3034019
$ws.Range("AT123").Value = 14.2592257
$ws.Range("AU123").Value = 20.4802587
$ws.Range("AV123").Value = 13.7841213
$ws.Range("AW123").Value = 14.4293391
$ws.Range("AX123").Value = 16.3772043
$ws.Range("AY123").Value = 15.2842832
$ws.Range("BA123").Value = 9.915507
$ws.Range("BB123").Value = 13.3723583
$ws.Range("BC123").Value = 14.4291145
$ws.Range("BD123").Value = 13.6497005
$ws.Range("BE123").Value = 15.5304343
# Row 124
$ws.Range("B124").Value = 12.3400853
$ws.Range("C124").Value = 18.6678701
$ws.Range("D124").Value = 17.4503912
$ws.Range("F124").Value = 15.1337443
$ws.Range("G124").Value = 11.5093612
$ws.Range("H124").Value = 15.2753969
$ws.Range("I124").Value = 15.5308322
$ws.Range("J124").Value = 16.277713
$ws.Range("K124").Value = 14.8768283
$ws.Range("L124").Value = 11.7416337
$ws.Range("M124").Value = 15.8710085
$ws.Range("O124").Value = 8.6385179
$ws.Range("P124").Value = 21.4670109
$ws.Range("Q124").Value = 13.1892586
$ws.Range("R124").Value = 16.0778909
$ws.Range("S124").Value = 18.7495091
$ws.Range("T124").Value = 15.0236843
$ws.Range("U124").Value = 15.0369839
$ws.Range("V124").Value = 19.3939806
$ws.Range("W124").Value = 16.8338189
$ws.Range("X124").Value = 15.8282473
$ws.Range("Y124").Value = 10.4820628
$ws.Range("Z124").Value = 14.0770636
$ws.Range("AA124").Value = 17.6398343
$ws.Range("AB124").Value = 14.9537741
$ws.Range("AD124").Value = 19.9505473
$ws.Range("AE124").Value = 10.2917635
$ws.Range("AF124").Value = 13.9127098
$ws.Range("AG124").Value = 17.5774844
$ws.Range("AH124").Value = 21.227876
$ws.Range("AI124").Value = 13.6594523
$ws.Range("AJ124").Value = 17.1309745
$ws.Range("AK124").Value = 13.9115842
$ws.Range("AL124").Value = 12.8987259
$ws.Range("AM124").Value = 15.2994503
$ws.Range("AN124").Value = 14.2817472
$ws.Range("AO124").Value = 13.9328542
$ws.Range("AP124").Value = 11.2304792
$ws.Range("AQ124").Value = 12.1814004
$ws.Range("AR124").Value = 12.2159577
$ws.Range("AS124").Value = 15.0243272
$ws.Range("AT124").Value = 14.3803623
$ws.Range("AU124").Value = 21.0610408
$ws.Range("AV124").Value = 13.6898051
$ws.Range("AW124").Value = 14.4886775
$ws.Range("AX124").Value = 16.6109142
$ws.Range("AY124").Value = 15.4748305
$ws.Range("BA124").Value = 9.5398067
$ws.Range("BB124").Value = 13.406725
$ws.Range("BC124").Value = 14.4690956
$ws.Range("BD124").Value = 13.7298209
$ws.Range("BE124").Value = 14.4995233
# Row 125
$ws.Range("B125").Value = 11.814693
$ws.Range("C125").Value = 18.7363156
$ws.Range("D125").Value = 17.670227
$ws.Range("F125").Value = 15.0548898
$ws.Range("G125").Value = 11.3314312
$ws.Range("H125").Value = 15.0712037
$ws.Range("I125").Value = 15.5639858
$ws.Range("J125").Value = 16.6083916
$ws.Range("K125").Value = 14.6426827
$ws.Range("L125").Value = 11.8639727
$ws.Range("M125").Value = 16.2637838
$ws.Range("O125").Value = 9.25
$ws.Range("P125").Value = 20.4053137
$ws.Range("Q125").Value = 13.3934708
$ws.Range("R125").Value = 15.8668252
$ws.Range("S125").Value = 18.6378651
$ws.Range("T125").Value = 14.8348661
$ws.Range("U125").Value = 14.7523711
$ws.Range("V125").Value = 19.5370663
$ws.Range("W125").Value = 16.0712625
$ws.Range("X125").Value = 16.2463057
$ws.Range("Y125").Value = 10.1835853
$ws.Range("Z125").Value = 14.0627452
$ws.Range("AA125").Value = 17.2956312
$ws.Range("AB125").Value = 14.8490425
$ws.Range("AD125").Value = 20.4356436
$ws.Range("AE125").Value = 9.4163891
$ws.Range("AF125").Value = 14.1992598
$ws.Range("AG125").Value = 18.5768029
$ws.Range("AH125").Value = 20.4326523
$ws.Range("AI125").Value = 13.4968638
$ws.Range("AJ125").Value = 16.9499312
$ws.Range("AK125").Value = 13.5927795
$ws.Range("AL125").Value = 12.5812421
$ws.Range("AM125").Value = 15.1868031
$ws.Range("AN125").Value = 14.0842797
$ws.Range("AO125").Value = 14.1265401
$ws.Range("AP125").Value = 11.4928076
$ws.Range("AQ125").Value = 12.0117758
$ws.Range("AR125").Value = 12.7329193
$ws.Range("AS125").Value = 14.347615
$ws.Range("AT125").Value = 14.3411821
$ws.Range("AU125").Value = 20.3420182
$ws.Range("AV125").Value = 13.784677
$ws.Range("AW125").Value = 14.5458575
$ws.Range("AX125").Value = 17.1316031
$ws.Range("AY125").Value = 15.1200427
$ws.Range("BA125").Value = 9.0402217
$ws.Range("BB125").Value = 13.4631368
$ws.Range("BC125").Value = 14.2372689
$ws.Range("BD125").Value = 13.7652473
$ws.Range("BE125").Value = 14.390123
# Row 126
$ws.Range("B126").Value = 11.3832853
$ws.Range("C126").Value = 19.5490362
$ws.Range("D126").Value = 17.8246431
$ws.Range("F126").Value = 15.1424916
$ws.Range("G126").Value = 11.3773956
$ws.Range("H126").Value = 14.6181609
$ws.Range("I126").Value = 15.3120165
$ws.Range("J126").Value = 15.1728553
$ws.Range("K126").Value = 15.043592
$ws.Range("L126").Value = 11.6157393
$ws.Range("M126").Value = 15.8155882
$ws.Range("O126").Value = 8.7894249
$ws.Range("P126").Value = 20.1245942
$ws.Range("Q126").Value = 13.8141082
$ws.Range("R126").Value = 15.5056831
$ws.Range("S126").Value = 18.5247995
$ws.Range("T126").Value = 14.0363711
$ws.Range("U126").Value = 15.4925417
$ws.Range("V126").Value = 19.2268127
$ws.Range("W126").Value = 16.1316809
$ws.Range("X126").Value = 15.9695245
$ws.Range("Y126").Value = 10.2288174
$ws.Range("Z126").Value = 13.7075491
$ws.Range("AA126").Value = 17.2119012
$ws.Range("AB126").Value = 14.9067591
$ws.Range("AD126").Value = 20.2638172
$ws.Range("AE126").Value = 9.6901589
$ws.Range("AF126").Value = 14.3056784
$ws.Range("AG126").Value = 18.1919969
$ws.Range("AH126").Value = 19.7387855
$ws.Range("AI126").Value = 13.6083181
$ws.Range("AJ126").Value = 16.3335512
$ws.Range("AK126").Value = 13.917699
$ws.Range("AL126").Value = 12.9687109
$ws.Range("AM126").Value = 14.7900317
$ws.Range("AN126").Value = 14.0925129
$ws.Range("AO126").Value = 14.2541212
$ws.Range("AP126").Value = 11.5370577
$ws.Range("AQ126").Value = 12.0422027
$ws.Range("AR126").Value = 12.7587101
$ws.Range("AS126").Value = 14.1460022
$ws.Range("AT126").Value = 14.4097329
$ws.Range("AU126").Value = 20.6840599
$ws.Range("AV126").Value = 13.6182758
$ws.Range("AW126").Value = 14.6861059
$ws.Range("AX126").Value = 17.4774932
$ws.Range("AY126").Value = 14.901903
$ws.Range("BA126").Value = 9.615384600000001
$ws.Range("BB126").Value = 13.2416785
$ws.Range("BC126").Value = 14.2862764
$ws.Range("BD126").Value = 12.8883871
$ws.Range("BE126").Value = 15.2005792
# Row 127
$ws.Range("B127").Value = 12.0035566
$ws.Range("C127").Value = 19.7044027
$ws.Range("D127").Value = 17.9351663
$ws.Range("F127").Value = 15.6301729
$ws.Range("G127").Value = 11.3313224
$ws.Range("H127").Value = 14.6786272
$ws.Range("I127").Value = 14.8741777
$ws.Range("J127").Value = 16.1782662
$ws.Range("K127").Value = 14.8640506
$ws.Range("L127").Value = 11.5611115
$ws.Range("M127").Value = 16.0075127
$ws.Range("O127").Value = 8.7438424
$ws.Range("P127").Value = 20.119043
$ws.Range("Q127").Value = 14.0976034
$ws.Range("R127").Value = 15.0949068
$ws.Range("S127").Value = 18.27988
$ws.Range("T127").Value = 14.365135
$ws.Range("U127").Value = 15.3878522
$ws.Range("V127").Value = 19.2828626
$ws.Range("W127").Value = 15.943804
$ws.Range("X127").Value = 15.8810978
$ws.Range("Y127").Value = 10.5403525
$ws.Range("Z127").Value = 13.5969602
$ws.Range("AA127").Value = 16.3080427
$ws.Range("AB127").Value = 14.8463573
$ws.Range("AD127").Value = 19.9377031
$ws.Range("AE127").Value = 9.8438313
$ws.Range("AF127").Value = 14.4325211
$ws.Range("AG127").Value = 18.0181321
$ws.Range("AH127").Value = 20.3556036
$ws.Range("AI127").Value = 12.6539624
$ws.Range("AJ127").Value = 15.9347273
$ws.Range("AK127").Value = 14.0268622
$ws.Range("AL127").Value = 12.7631282
$ws.Range("AM127").Value = 14.4957168
$ws.Range("AN127").Value = 13.72005
$ws.Range("AO127").Value = 14.4480335
$ws.Range("AP127").Value = 11.4184249
$ws.Range("AQ127").Value = 11.8629735
$ws.Range("AR127").Value = 10.3846154
$ws.Range("AS127").Value = 13.4677713
$ws.Range("AT127").Value = 14.319057
$ws.Range("AU127").Value = 19.252603
$ws.Range("AV127").Value = 13.4383922
$ws.Range("AW127").Value = 14.5067304
$ws.Range("AX127").Value = 17.6116455
$ws.Range("AY127").Value = 14.922335
$ws.Range("BA127").Value = 8.993971999999999
$ws.Range("BB127").Value = 12.9343328
$ws.Range("BC127").Value = 13.9741451
$ws.Range("BD127").Value = 13.6967255
$ws.Range("BE127").Value = 15.7766241
# Row 128
$ws.Range("B128").Value = 13.0201765
$ws.Range("C128").Value = 19.3907246
$ws.Range("D128").Value = 18.3076257
$ws.Range("F128").Value = 15.7920972
$ws.Range("G128").Value = 11.4733163
$ws.Range("H128").Value = 14.5098615
$ws.Range("I128").Value = 14.4362663
$ws.Range("J128").Value = 16.2077597
$ws.Range("K128").Value = 14.7579948
$ws.Range("L128").Value = 11.774065
$ws.Range("M128").Value = 15.7653073
$ws.Range("O128").Value = 7.485349
$ws.Range("P128").Value = 20.3072495
$ws.Range("Q128").Value = 13.5186865
$ws.Range("R128").Value = 15.1311186
$ws.Range("S128").Value = 18.606122
$ws.Range("T128").Value = 14.5635459
$ws.Range("U128").Value = 15.2259707
$ws.Range("V128").Value = 19.2080921
$ws.Range("W128").Value = 16.1783622
$ws.Range("X128").Value = 15.674513
$ws.Range("Y128").Value = 10.5034265
$ws.Range("Z128").Value = 13.4264022
$ws.Range("AA128").Value = 16.7609852
$ws.Range("AB128").Value = 14.9009109
$ws.Range("AD128").Value = 20.8144855
$ws.Range("AE128").Value = 10.4296009
$ws.Range("AF128").Value = 14.3535007
$ws.Range("AG128").Value = 19.6169035
$ws.Range("AH128").Value = 20.3654683
$ws.Range("AI128").Value = 13.0580958
$ws.Range("AJ128").Value = 15.4762853
$ws.Range("AK128").Value = 13.2180148
$ws.Range("AL128").Value = 12.9405504
$ws.Range("AM128").Value = 14.5441568
$ws.Range("AN128").Value = 13.5761371
$ws.Range("AO128").Value = 13.8899385
$ws.Range("AP128").Value = 11.1351554
$ws.Range("AQ128").Value = 11.7157836
$ws.Range("AS128").Value = 13.9287102
$ws.Range("AT128").Value = 14.4988624
$ws.Range("AU128").Value = 19.9558536
$ws.Range("AV128").Value = 13.7041421
$ws.Range("AW128").Value = 14.4571091
$ws.Range("AX128").Value = 17.8783984
$ws.Range("AY128").Value = 14.3251903
$ws.Range("BA128").Value = 9.468006900000001
$ws.Range("BB128").Value = 12.8820076
$ws.Range("BC128").Value = 13.976316
$ws.Range("BD128").Value = 13.5962756
$ws.Range("BE128").Value = 15.8181815
